$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated values in column F (dSF) for the given rows, per repull/push of data
$updates = @{
    2  = 0
    5  = 0
    13 = -2
    15 = -1
    16 = -4
    17 = -2
    27 = -2
    28 = -4
    31 = 0
    32 = -4
    36 = -2
    43 = -2
    47 = 0
    52 = -3
    54 = -1
    57 = -8
    58 = 6
    61 = -2
    64 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
